$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 990
$ws.Range("I12").Value = 990
$ws.Range("K12").Value = 990
$ws.Range("M12").Value = -820
$ws.Range("H17").Value = 669.4666999999999
$ws.Range("J17").Value = 669.4666999999999
$ws.Range("L17").Value = 2008.4001
$ws.Range("N17").Value = -2344.4001
$ws.Range("H40").Value = 5002
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 5002
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 5002
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -5352
$ws.Range("H101").Value = 646
$ws.Range("I101").Value = 654
$ws.Range("K101").Value = 1962
$ws.Range("M101").Value = -340
$ws.Range("H138").Value = 5966.4375
$ws.Range("J138").Value = 6738.475
$ws.Range("L138").Value = 20215.425
$ws.Range("N138").Value = -30495.425

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4849.4
$ws.Range("I61").Value = 5054.8887
$ws.Range("K61").Value = 5054.8887
$ws.Range("M61").Value = -4842.8887
$ws.Range("H74").Value = 874.2727
$ws.Range("I74").Value = 846.44446
$ws.Range("K74").Value = 846.44446
$ws.Range("M74").Value = 27.55553999999995
$ws.Range("H77").Value = 874.2727
$ws.Range("I77").Value = 846.44446
$ws.Range("K77").Value = 4232.2223
$ws.Range("M77").Value = 135.7776999999996
$ws.Range("H122").Value = 4573.769
$ws.Range("I122").Value = 4557.56
$ws.Range("K122").Value = 13672.68
$ws.Range("M122").Value = -11222.68
$ws.Range("H132").Value = 3856.1428
$ws.Range("I132").Value = 3165.3333
$ws.Range("J132").Value = 4374.25
$ws.Range("K132").Value = 9495.999899999999
$ws.Range("L132").Value = 13122.75
$ws.Range("M132").Value = -6965.999899999999
$ws.Range("N132").Value = -18182.75
$ws.Range("H136").Value = 4849.4
$ws.Range("I136").Value = 5054.8887
$ws.Range("K136").Value = 15164.6661
$ws.Range("M136").Value = -12614.6661

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H134").Value = 4854.5454
$ws.Range("J134").Value = 4499.5
$ws.Range("L134").Value = 13498.5
$ws.Range("N134").Value = -18568.5

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 12736.75
$ws.Range("J43").Value = 12736.75
$ws.Range("L43").Value = 12736.75
$ws.Range("N43").Value = -13104.75
$ws.Range("H58").Value = 792.5
$ws.Range("I58").Value = 792.5
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 792.5
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -589.5
$ws.Range("N58").ClearContents()
$ws.Range("H74").Value = 34469.75
$ws.Range("J74").Value = 34469.75
$ws.Range("L74").Value = 34469.75
$ws.Range("N74").Value = -36217.75
$ws.Range("H77").Value = 34469.75
$ws.Range("J77").Value = 34469.75
$ws.Range("L77").Value = 103409.25
$ws.Range("N77").Value = -112145.25
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("H101").Value = 12736.75
$ws.Range("J101").Value = 12736.75
$ws.Range("L101").Value = 12736.75
$ws.Range("N101").Value = -19226.75
$ws.Range("H132").Value = 4398.222
$ws.Range("I132").Value = 2366.3333
$ws.Range("J132").Value = 5414.1665
$ws.Range("K132").Value = 7098.999899999999
$ws.Range("L132").Value = 16242.4995
$ws.Range("M132").Value = -4568.999899999999
$ws.Range("N132").Value = -21302.4995
$ws.Range("H134").Value = 4491.0586
$ws.Range("I134").Value = 4390.1333
$ws.Range("K134").Value = 13170.3999
$ws.Range("M134").Value = -10635.3999
$ws.Range("H136").Value = 792.5
$ws.Range("I136").Value = 792.5
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 2377.5
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = 172.5
$ws.Range("N136").ClearContents()

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 5000
$ws.Range("I69").Value = 4000
$ws.Range("J69").Value = 6000
$ws.Range("K69").Value = 12000
$ws.Range("L69").Value = 18000
$ws.Range("M69").Value = -11189
$ws.Range("N69").Value = -19622
$ws.Range("H72").Value = 5000
$ws.Range("I72").Value = 4000
$ws.Range("J72").Value = 6000
$ws.Range("K72").Value = 36000
$ws.Range("L72").Value = 54000
$ws.Range("M72").Value = -31944
$ws.Range("N72").Value = -62112
$ws.Range("H113").Value = 770.6667
$ws.Range("I113").Value = 649.5
$ws.Range("J113").Value = 831.25
$ws.Range("K113").Value = 1948.5
$ws.Range("L113").Value = 2493.75
$ws.Range("M113").Value = 221.5
$ws.Range("N113").Value = -6833.75

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 718
$ws.Range("I107").Value = 1000
$ws.Range("J107").Value = 436
$ws.Range("K107").Value = 1000
$ws.Range("L107").Value = 436
$ws.Range("M107").Value = 920
$ws.Range("N107").Value = -4276
$ws.Range("H113").Value = 2490
$ws.Range("I113").Value = 2490
$ws.Range("K113").Value = 2490
$ws.Range("M113").Value = -320
$ws.Range("H132").Value = 7849.1665
$ws.Range("I132").Value = 8774.25
$ws.Range("K132").Value = 26322.75
$ws.Range("M132").Value = -23792.75

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4990
$ws.Range("I46").Value = 4990
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 4990
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -4802
$ws.Range("N46").ClearContents()
$ws.Range("H61").Value = 2291.7856
$ws.Range("I61").Value = 2214.6155
$ws.Range("K61").Value = 2214.6155
$ws.Range("M61").Value = -2012.6155
$ws.Range("H113").Value = 2291.7856
$ws.Range("I113").Value = 2214.6155
$ws.Range("K113").Value = 2214.6155
$ws.Range("M113").Value = -44.61549999999988
$ws.Range("H136").Value = 2666.3333
$ws.Range("I136").Value = 1999.5
$ws.Range("K136").Value = 5998.5
$ws.Range("M136").Value = -3448.5

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 100370
$ws.Range("J104").Value = 100370
$ws.Range("L104").Value = 100370
$ws.Range("N104").Value = -107358
$ws.Range("H113").Value = 1333
$ws.Range("I113").Value = 999
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 2997
$ws.Range("L113").Value = 4500
$ws.Range("M113").Value = -827
$ws.Range("N113").Value = -8840
$ws.Range("H132").Value = 3379.9375
$ws.Range("I132").Value = 2204.2856
$ws.Range("J132").Value = 4294.3335
$ws.Range("K132").Value = 6612.8568
$ws.Range("L132").Value = 12883.0005
$ws.Range("M132").Value = -4082.8568
$ws.Range("N132").Value = -17943.0005
$ws.Range("H136").Value = 7980.2666
$ws.Range("I136").Value = 7927.037
$ws.Range("K136").Value = 23781.111
$ws.Range("M136").Value = -21231.111
